$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "65.090.69"
Set-TextValue $ws.Range("E2") "  +0.98%  "
Set-TextValue $ws.Range("D3") "3.378.47"
Set-TextValue $ws.Range("E3") "  +0.61%  "
Set-TextValue $ws.Range("E4") "  +0.00%  "
Set-TextValue $ws.Range("D5") "555.27"
Set-TextValue $ws.Range("E5") "  -0.12%  "
Set-TextValue $ws.Range("D6") "174.10"
Set-TextValue $ws.Range("E6") "  -0.74%  "
Set-TextValue $ws.Range("E7") "  +2.09%  "
Set-TextValue $ws.Range("D8") "3.368.32"
Set-TextValue $ws.Range("E8") "  +0.58%  "
Set-TextValue $ws.Range("E9") "  -0.01%  "
Set-TextValue $ws.Range("D10") "0.173"
Set-TextValue $ws.Range("E10") "  +5.66%  "
Set-TextValue $ws.Range("D11") "0.636"
Set-TextValue $ws.Range("E11") "  +1.18%  "
Set-TextValue $ws.Range("D12") "53.61"
Set-TextValue $ws.Range("E12") "  -1.59%  "
Set-TextValue $ws.Range("E13") "  +1.13%  "
Set-TextValue $ws.Range("D14") "9.15"
Set-TextValue $ws.Range("E14") "  +0.64%  "
Set-TextValue $ws.Range("D15") "3.917.68"
Set-TextValue $ws.Range("E15") "  +0.76%  "
Set-TextValue $ws.Range("D16") "18.32"
Set-TextValue $ws.Range("E16") "  -0.55%  "
Set-TextValue $ws.Range("B17") "WrappedEther"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "3.398.78"
Set-TextValue $ws.Range("E17") "  +1.26%  "
Set-TextValue $ws.Range("B18") "TRON"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D18") "0.118"
Set-TextValue $ws.Range("E18") "  +0.12%  "
Set-TextValue $ws.Range("D19") "65.002.19"
Set-TextValue $ws.Range("E19") "  +0.94%  "
Set-TextValue $ws.Range("D20") "11.82"
Set-TextValue $ws.Range("E20") "  -0.23%  "
Set-TextValue $ws.Range("E21") "  +1.12%  "
Set-TextValue $ws.Range("D22") "455.65"
Set-TextValue $ws.Range("E22") "  -1.35%  "
Set-TextValue $ws.Range("D23") "4.87"
Set-TextValue $ws.Range("E23") "  +0.44%  "
Set-TextValue $ws.Range("E24") "  -0.44%  "
Set-TextValue $ws.Range("E25") "  +6.14%  "
Set-TextValue $ws.Range("D26") "87.65"
Set-TextValue $ws.Range("E26") "  +1.69%  "
Set-TextValue $ws.Range("E27") "  +1.25%  "
Set-TextValue $ws.Range("D28") "10.67"
Set-TextValue $ws.Range("E28") "  -2.61%  "
Set-TextValue $ws.Range("D29") "8.69"
Set-TextValue $ws.Range("E29") "  -1.08%  "
Set-TextValue $ws.Range("D30") "31.10"
Set-TextValue $ws.Range("E30") "  +3.31%  "
Set-TextValue $ws.Range("E31") "  -2.04%  "
Set-TextValue $ws.Range("D32") "63.30"
Set-TextValue $ws.Range("E32") "  +7.62%  "
Set-TextValue $ws.Range("D33") "11.45"
Set-TextValue $ws.Range("E33") "  -0.23%  "
Set-TextValue $ws.Range("D34") "577.84"
Set-TextValue $ws.Range("E34") "  -0.69%  "
Set-TextValue $ws.Range("E35") "  -0.71%  "
Set-TextValue $ws.Range("E36") "  -0.03%  "
Set-TextValue $ws.Range("D37") "3.60"
Set-TextValue $ws.Range("E37") "  +2.57%  "
Set-TextValue $ws.Range("D38") "0.142"
Set-TextValue $ws.Range("E38") "  +1.69%  "
Set-TextValue $ws.Range("D39") "35.65"
Set-TextValue $ws.Range("E39") "  -0.05%  "
Set-TextValue $ws.Range("E40") "  -0.74%  "
Set-TextValue $ws.Range("E41") "  -2.68%  "
Set-TextValue $ws.Range("D42") "3.104.76"
Set-TextValue $ws.Range("E42") "  +0.26%  "
Set-TextValue $ws.Range("E43") "  +1.14%  "
Set-TextValue $ws.Range("D44") "2.76"
Set-TextValue $ws.Range("E44") "  -1.44%  "
Set-TextValue $ws.Range("E45") "  -0.48%  "
Set-TextValue $ws.Range("E46") "  +2.15%  "
Set-TextValue $ws.Range("E47") "  -3.29%  "
Set-TextValue $ws.Range("E48") "  -0.05%  "
Set-TextValue $ws.Range("D49") "140.57"
Set-TextValue $ws.Range("E49") "  +3.83%  "
Set-TextValue $ws.Range("E50") "  -2.36%  "
Set-TextValue $ws.Range("D51") "8.31"
Set-TextValue $ws.Range("E51") "  -1.00%  "
